$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.093714853270462584, 0.093531712790372978),
    @(-0.084097524412325875, 0.083683130729807154),
    @(-0.033977177575641448, 0.033951005705759485),
    @(-0.02595100583330634, 0.025828150869527278),
    @(-0.022828150929305124, 0.022429891428383364),
    @(-0.0086760630258648774, 0.0085270224705027431),
    @(0.0014729773645716371, -0.0014972255512262223),
    @(0.011497225386481347, -0.011523252480719393),
    @(0.013523252431637101, -0.013541610853327057),
    @(0.015541610807703776, -0.01554100317556717),
    @(0.018541003115952748, -0.018541766544696081),
    @(0.022041766479957925, -0.022078446983656352),
    @(0.017846492813091253, -0.017986732195576671),
    @(0.025986732078767893, -0.02606369135902753),
    @(-0.0080494065617493504, 0.0080324206535324905),
    @(-0.0060324206831596783, 0.0060026646607012069),
    @(-0.0040026646920985343, 0.0039999999400555097),
    @(-0.0034341884640234355, 0.003360626855187121),
    @(0.0006393730867793046, -0.0012011640716531602),
    @(0.005201164012420989, -0.0053562533537920842),
    @(0.0093562532948716637, -0.0095951682939476868),
    @(-0.04570595322301152, 0.045495200521422419),
    @(-0.040495200609847792, 0.040097956956954839),
    @(-0.020097957267005917, 0.019999999685413883),
    @(-0.011753828454326154, 0.011753862112023228),
    @(-0.0092538621711621971, 0.0092535631601275981),
    @(-0.0067535632196875639, 0.0067423213653805547),
    @(-0.0047423214201254282, 0.0047402808146204123),
    @(0.0022597190591087113, -0.0022595366654538651),
    @(0.062259535792228871, -0.062542993638926347),
    @(0.046068838155784775, -0.046145434256697371),
    @(-0.0040007170905571598, 0.0039999999357274163)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
